# "corrigido bug de celulas vazias" - fix for the empty-cells bug in the
# Contatos sheet: refresh the contact rows and make sure row 3 carries
# formatted (if empty) trailing cells out to column E, so downstream
# readers that expect a rectangular A1:E3 range no longer choke on a
# short/ragged row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the used range with three formatted-but-empty columns (C:E).
# Inserting columns copies row 3's existing formatting into the new
# cells, which is exactly the "missing trailing cells" fix; rows 1-2
# should stay untouched, so their spilled-over inserted cells are
# cleared right back out.
$ws.Columns("C:E").Insert()
$ws.Range("C1:E2").Clear()

# Replace the old demo/test contacts with the current ones.
$ws.Range("A2").Value = "Andre MFprint "
$ws.Range("B2").Value = "  ola"
$ws.Range("A3").Value = "Nathaly"
$ws.Range("B3").Value = "ola teste"

# Widen the message column slightly so the new text isn't clipped.
$ws.Columns("B").ColumnWidth = 10.43
